$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.040.30"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "2.952.48"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'378.52"
$ws.Range("D6").Value = "'101.59"
$ws.Range("E6").Value = "  -2.45%  "
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.585"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "3.415.00"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").Value = "'7.37"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "2.941.99"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  +5.24%  "
$ws.Range("D18").Value = "51.113.07"
$ws.Range("D19").Value = "'3.21"
$ws.Range("E19").Value = "  -6.51%  "
$ws.Range("D20").Value = "'7.16"
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("E21").Value = "  -3.81%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "'68.41"
$ws.Range("D24").Value = "'261.02"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("E25").Value = "  +2.98%  "
$ws.Range("D26").Value = "'8.25"
$ws.Range("E26").Value = "  +11.92%  "
$ws.Range("D27").Value = "'7.71"
$ws.Range("E27").Value = "  +9.61%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'4.10"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").Value = "'0.112"
$ws.Range("E31").Value = "  +9.92%  "
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").Value = "'9.77"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").Value = "'50.52"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("D36").Value = "'33.50"
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("E37").Value = "  +2.64%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("D40").Value = "'16.82"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").Value = "'121.49"
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("D45").Value = "'20.99"
$ws.Range("E45").Value = "  -3.91%  "
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("D47").Value = "'0.274"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D49").Value = "2.003.38"
$ws.Range("E49").Value = "  -1.19%  "
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("E51").Value = "  +4.41%  "
